$wb = $excel.ActiveWorkbook

# --- Metadata sheet: URL / Version / Date / Publisher updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/admit-count-epis"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet: clear the stray Constraint(s) text duplicated on the
#     "Extension" row (it belongs only on the "Extension.extension" row) ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""

# The "Extension.url" row's Fixed Value mirrors the StructureDefinition's own
# URL, so it must be kept in sync with the Metadata sheet's URL update above.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/admit-count-epis"
